$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 698, shifting existing rows 698-762 down to 699-763
$ws.Rows.Item(698).Insert()

# Populate the new row 698 with the new market price record
$ws.Cells.Item(698, 1).Value = 10
$ws.Cells.Item(698, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(698, 3).Value = "La Araucanía"
$ws.Cells.Item(698, 4).Value = 45223
$ws.Cells.Item(698, 5).Value = 9
$ws.Cells.Item(698, 6).Value = 100112028
$ws.Cells.Item(698, 7).Value = "Sandia"
$ws.Cells.Item(698, 8).Value = "Sin especificar"
$ws.Cells.Item(698, 9).Value = "Primera"
$ws.Cells.Item(698, 10).Value = 180
$ws.Cells.Item(698, 11).Value = 3000
$ws.Cells.Item(698, 12).Value = 4000
$ws.Cells.Item(698, 13).Value = 3556
$ws.Cells.Item(698, 14).Value = "$/unidad"
$ws.Cells.Item(698, 15).Value = "Brasil"
$ws.Cells.Item(698, 16).Value = 3556
$ws.Cells.Item(698, 17).Value = 1
$ws.Cells.Item(698, 18).Value = "Hortaliza"
